$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell value for A3
$ws.Range("A3").Value = 2

# Select A3 to match the saved selection/active cell in the diff
$ws.Range("A3").Select()
